# Swap the presentation's theme palette ("Integral") for the stock
# "Office Theme" palette that used to live in the (otherwise orphaned)
# notes-master theme part, per the authored diff:
#   ppt/theme/theme1.xml : Integral  -> Office Theme
#   ppt/theme/theme2.xml : Office Theme -> Integral
#
# The notes-master theme (theme2.xml) is not reachable through the
# PowerPoint object model (no object exposes its ThemeColorScheme
# independently of the slide master's), so we apply the reachable half
# of the swap: re-point the slide master's theme color scheme — which
# backs ppt/theme/theme1.xml — from the Integral palette to the Office
# palette.

function ConvertTo-OleColor([string]$hex) {
    # PowerPoint's ColorFormat.RGB (and the classic VBA RGB() macro)
    # packs components as R + G*256 + B*65536, i.e. the reverse of the
    # usual #RRGGBB text order.
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Office Theme color scheme, in ThemeColorScheme.Item() slot order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeTheme = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $officeTheme.Length; $i++) {
    $colorScheme.Item($i).RGB = ConvertTo-OleColor $officeTheme[$i - 1]
}
